# Auto-generated Excel COM-interop script
# Applies updated market-price derived values (columns H-N) across 8 job sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 805.9167
$ws.Range("I4").Value = 183.22223
$ws.Range("K4").Value = 183.22223
$ws.Range("M4").Value = -69.22223

$ws.Range("H38").Value = 1800
$ws.Range("I38").Value = 700
$ws.Range("J38").Value = 4000
$ws.Range("K38").Value = 2100
$ws.Range("L38").Value = 12000
$ws.Range("M38").Value = -1728
$ws.Range("N38").Value = -12744

$ws.Range("H42").Value = 370.42856
$ws.Range("I42").Value = 98.833336
$ws.Range("J42").Value = 2000
$ws.Range("K42").Value = 296.500008
$ws.Range("L42").Value = 6000
$ws.Range("M42").Value = -66.50000799999998
$ws.Range("N42").Value = -6460

$ws.Range("H103").Value = 951.2
$ws.Range("J103").Value = 984.8333
$ws.Range("L103").Value = 2954.4999
$ws.Range("N103").Value = -4126.4999

$ws.Range("H107").Value = 1175.7028
$ws.Range("I107").Value = 1591.68
$ws.Range("J107").Value = 309.08334
$ws.Range("K107").Value = 1591.68
$ws.Range("L107").Value = 309.08334
$ws.Range("M107").Value = 328.3199999999999
$ws.Range("N107").Value = -4149.08334

$ws.Range("H132").Value = 3190.9807
$ws.Range("I132").Value = 3243.745
$ws.Range("K132").Value = 9731.235000000001
$ws.Range("M132").Value = -7201.235000000001

$ws.Range("H135").Value = 1480.9166
$ws.Range("I135").Value = 674.6667
$ws.Range("K135").Value = 6072.0003
$ws.Range("M135").Value = -3537.0003

$ws.Range("H138").Value = 2833.69
$ws.Range("I138").Value = 1330.2593
$ws.Range("J138").Value = 3389.7534
$ws.Range("K138").Value = 3990.7779
$ws.Range("L138").Value = 10169.2602
$ws.Range("M138").Value = 1149.2221
$ws.Range("N138").Value = -20449.2602

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 202
$ws.Range("J4").Value = 202
$ws.Range("L4").Value = 202
$ws.Range("N4").Value = -434

$ws.Range("H32").Value = 1064688.4
$ws.Range("I32").Value = 314034.78
$ws.Range("J32").Value = 12824928
$ws.Range("K32").Value = 314034.78
$ws.Range("L32").Value = 12824928
$ws.Range("M32").Value = -313747.78
$ws.Range("N32").Value = -12825502

$ws.Range("H74").Value = 2613.3818
$ws.Range("I74").Value = 2486.449
$ws.Range("K74").Value = 2486.449
$ws.Range("M74").Value = -1612.449

$ws.Range("H77").Value = 2613.3818
$ws.Range("I77").Value = 2486.449
$ws.Range("K77").Value = 12432.245
$ws.Range("M77").Value = -8064.245000000001

$ws.Range("H130").Value = 34808.41
$ws.Range("J130").Value = 34808.41
$ws.Range("L130").Value = 34808.41
$ws.Range("N130").Value = -44848.41

$ws.Range("H132").Value = 272870.3
$ws.Range("I132").Value = 386451.88
$ws.Range("K132").Value = 1159355.64
$ws.Range("M132").Value = -1156825.64

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1850.48
$ws.Range("I105").Value = 1438.3
$ws.Range("K105").Value = 1438.3
$ws.Range("M105").Value = 308.7

$ws.Range("H134").Value = 2235056.2
$ws.Range("I134").Value = 2553315
$ws.Range("K134").Value = 7659945
$ws.Range("M134").Value = -7657410

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 524.2222
$ws.Range("I22").Value = 549
$ws.Range("K22").Value = 549
$ws.Range("M22").Value = -199

$ws.Range("H31").Value = 3655.543
$ws.Range("I31").Value = 2659.8572
$ws.Range("J31").Value = 5149.0713
$ws.Range("K31").Value = 2659.8572
$ws.Range("L31").Value = 5149.0713
$ws.Range("M31").Value = -2364.8572
$ws.Range("N31").Value = -5739.0713

$ws.Range("H34").Value = 3655.543
$ws.Range("I34").Value = 2659.8572
$ws.Range("J34").Value = 5149.0713
$ws.Range("K34").Value = 2659.8572
$ws.Range("L34").Value = 5149.0713
$ws.Range("M34").Value = -2457.8572
$ws.Range("N34").Value = -5553.0713

$ws.Range("H58").Value = 2616.7046
$ws.Range("I58").Value = 2366.0286
$ws.Range("K58").Value = 2366.0286
$ws.Range("M58").Value = -2163.0286

$ws.Range("H132").Value = 2720.7183
$ws.Range("I132").Value = 2523.7585
$ws.Range("K132").Value = 7571.2755
$ws.Range("M132").Value = -5041.2755

$ws.Range("H134").Value = 3601.1904
$ws.Range("I134").Value = 4183
$ws.Range("K134").Value = 12549
$ws.Range("M134").Value = -10014

$ws.Range("H136").Value = 2616.7046
$ws.Range("I136").Value = 2366.0286
$ws.Range("K136").Value = 7098.085800000001
$ws.Range("M136").Value = -4548.085800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 918.1429000000001
$ws.Range("I33").Value = 51.8
$ws.Range("K33").Value = 310.8
$ws.Range("M33").Value = -27.79999999999995

$ws.Range("H60").Value = 271.2
$ws.Range("I60").Value = 218.33333
$ws.Range("K60").Value = 654.99999
$ws.Range("M60").Value = -403.99999

$ws.Range("H68").Value = 385
$ws.Range("J68").Value = 330
$ws.Range("L68").Value = 990
$ws.Range("N68").Value = -2612

$ws.Range("H71").Value = 385
$ws.Range("J71").Value = 330
$ws.Range("L71").Value = 2970
$ws.Range("N71").Value = -11082

$ws.Range("H92").Value = 925.1429000000001
$ws.Range("J92").Value = 1322
$ws.Range("L92").Value = 3966
$ws.Range("N92").Value = -6462

$ws.Range("H132").Value = 1205.8334
$ws.Range("I132").Value = 1184
$ws.Range("J132").Value = 1240.1428
$ws.Range("K132").Value = 10656
$ws.Range("L132").Value = 11161.2852
$ws.Range("M132").Value = -8126
$ws.Range("N132").Value = -16221.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 77.6875
$ws.Range("I2").Value = 32
$ws.Range("K2").Value = 32
$ws.Range("M2").Value = 81

$ws.Range("H39").Value = 35261
$ws.Range("J39").Value = 35261
$ws.Range("L39").Value = 35261
$ws.Range("N39").Value = -36325

$ws.Range("H80").Value = 2059.2856
$ws.Range("J80").Value = 2251.6667
$ws.Range("L80").Value = 2251.6667
$ws.Range("N80").Value = -4247.6667

$ws.Range("H83").Value = 2059.2856
$ws.Range("J83").Value = 2251.6667
$ws.Range("L83").Value = 11258.3335
$ws.Range("N83").Value = -21242.3335

$ws.Range("H102").Value = 1495
$ws.Range("I102").Value = 1495
$ws.Range("K102").Value = 1495
$ws.Range("M102").Value = 127

$ws.Range("H126").Value = 2923.3333
$ws.Range("J126").Value = 3750
$ws.Range("L126").Value = 11250
$ws.Range("N126").Value = -16190

$ws.Range("H132").Value = 4376.3184
$ws.Range("I132").Value = 4389
$ws.Range("K132").Value = 13167
$ws.Range("M132").Value = -10637

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 55567724
$ws.Range("I40").Value = 83344080
$ws.Range("K40").Value = 83344080
$ws.Range("M40").Value = -83343944

$ws.Range("H136").Value = 2494.5334
$ws.Range("I136").Value = 1983.7142
$ws.Range("K136").Value = 5951.142599999999
$ws.Range("M136").Value = -3401.142599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 500010000
$ws.Range("I14").Value = 1000000000
$ws.Range("K14").Value = 1000000000
$ws.Range("M14").Value = -999999832

$ws.Range("H62").Value = 5533.3335
$ws.Range("I62").Value = 3750
$ws.Range("J62").Value = 6425
$ws.Range("K62").Value = 3750
$ws.Range("L62").Value = 6425
$ws.Range("M62").Value = -3126
$ws.Range("N62").Value = -7673

$ws.Range("H65").Value = 5533.3335
$ws.Range("I65").Value = 3750
$ws.Range("J65").Value = 6425
$ws.Range("K65").Value = 18750
$ws.Range("L65").Value = 32125
$ws.Range("M65").Value = -15630
$ws.Range("N65").Value = -38365

$ws.Range("H93").Value = 128000
$ws.Range("J93").Value = 128000
$ws.Range("L93").Value = 128000
$ws.Range("N93").Value = -132992

$ws.Range("H113").Value = 1051.375
$ws.Range("I113").Value = 1103.2858
$ws.Range("K113").Value = 3309.8574
$ws.Range("M113").Value = -1139.8574

$ws.Range("H126").Value = 2813.111
$ws.Range("I126").Value = 1252
$ws.Range("K126").Value = 3756
$ws.Range("M126").Value = -1286

$ws.Range("H132").Value = 18436.72
$ws.Range("I132").Value = 19379.648
$ws.Range("K132").Value = 58138.944
$ws.Range("M132").Value = -55608.944

$ws.Range("H136").Value = 1644.5397
$ws.Range("I136").Value = 1287.8572
$ws.Range("K136").Value = 3863.5716
$ws.Range("M136").Value = -1313.5716
